$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Ost (Q) / Nord (R) coordinate values to whole numbers
$ws.Range("Q2").Value = 358653
$ws.Range("R2").Value = 6874558
$ws.Range("Q3").Value = 358671
$ws.Range("R3").Value = 6874549

# Clear the Starttid (Z) and Sluttid (AB) cells for rows 2 and 3
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
